# Auto-generated edit script for CATALOG_pics.xlsx "Files" sheet update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Files")

# --- Header row (row 1): insert new Thumbnail column (B), shift others, add TYPE: and JSON:exifData columns ---
$ws.Range("A1").Value = "FILE:Filename"
$ws.Range("B1").Value = "RELATION:Thumbnail"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "RELATION:Creator*"
$ws.Range("E1").Value = "RELATION:License"
$ws.Range("F1").Value = "RELATION:contentLocation"
$ws.Range("G1").Value = "TYPE:"
$ws.Range("H1").Value = "JSON:exifData"

# --- Row 2: update existing file record, add thumbnail path, new description, drop old License/Contributor/contentLocation values, add TYPE + exif JSON ---
$ws.Range("A2").Value = "2017-06-11 12.56.14.jpg"
$ws.Range("B2").Value = "pics/thumbs/2017-06-11 12.56.14.png"
$ws.Range("C2").Value = 'Depicts a fence at a disused motor racing venue with the front part of a slightly out of focus black dog in the foreground.'
$ws.Range("D2").Value = "Peter Sefton"
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = "ImageObject"
$ws.Range("H2").Value = '[{"@type": "PropertyValue", "value": "Off", "name": "AELock"}, {"@type": "PropertyValue", "value": "(111,106)-(144,149)", "name": "AFAreas"}, {"@type": "PropertyValue", "value": "Off", "name": "AFFineTune"}, {"@type": "PropertyValue", "value": "0 0 0", "name": "AFFineTuneAdj"}, {"@type": "PropertyValue", "value": "Left (or n/a)", "name": "AFPoint"}, {"@type": "PropertyValue", "value": "(49%,69%) (49%,69%)", "name": "AFPointSelected"}, {"@type": "PropertyValue", "value": "Ready", "name": "AFSearch"}, {"@type": "PropertyValue", "value": 6.3, "name": "Aperture"}, {"@type": "PropertyValue", "value": "Off; 0; 0; 0", "name": "ArtFilter"}, {"@type": "PropertyValue", "value": "0 0 4031 2687", "name": "AspectFrame"}, {"@type": "PropertyValue", "value": "3:2", "name": "AspectRatio"}, {"@type": "PropertyValue", "value": 8, "name": "BitsPerSample"}, {"@type": "PropertyValue", "value": "64 64 64 64", "name": "BlackLevel2"}, {"@type": "PropertyValue", "value": 1.367188, "name": "BlueBalance"}, {"@type": "PropertyValue", "value": 1.103, "name": "BodyFirmwareVersion"}, {"@type": "PropertyValue", "value": "OLYMPUS DIGITAL CAMERA         ", "name": "CameraID"}, {"@type": "PropertyValue", "value": "0100", "name": "CameraSettingsVersion"}, {"@type": "PropertyValue", "value": "E-PL1", "name": "CameraType2"}, {"@type": "PropertyValue", "value": "0.015 mm", "name": "CircleOfConfusion"}, {"@type": "PropertyValue", "value": 3, "name": "ColorComponents"}, {"@type": "PropertyValue", "value": "322 -40 -26 -32 332 -44 -2 -72 330", "name": "ColorMatrix"}, {"@type": "PropertyValue", "value": "sRGB", "name": "ColorSpace"}, {"@type": "PropertyValue", "value": "Y, Cb, Cr, -", "name": "ComponentsConfiguration"}, {"@type": "PropertyValue", "value": "JPEG (old-style)", "name": "Compression"}, {"@type": "PropertyValue", "value": 4, "name": "CompressionFactor"}, {"@type": "PropertyValue", "value": "Normal", "name": "Contrast"}, {"@type": "PropertyValue", "value": "0 (min -5, max 5)", "name": "ContrastSetting"}, {"@type": "PropertyValue", "value": "2017:06:11 12:56:14", "name": "CreateDate"}, {"@type": "PropertyValue", "value": 3024, "name": "CropHeight"}, {"@type": "PropertyValue", "value": "28 0", "name": "CropLeft"}, {"@type": "PropertyValue", "value": "30 0", "name": "CropTop"}, {"@type": "PropertyValue", "value": 4032, "name": "CropWidth"}, {"@type": "PropertyValue", "value": "Normal", "name": "CustomRendered"}, {"@type": "PropertyValue", "value": "1 (min -5, max 5)", "name": "CustomSaturation"}, {"@type": "PropertyValue", "value": "0.98 m (1.04 - 2.02 m)", "name": "DOF"}, {"@type": "PropertyValue", "value": "2017:06:11 12:56:14", "name": "DateTimeOriginal"}, {"@type": "PropertyValue", "value": 1, "name": "DigitalZoomRatio"}, {"@type": "PropertyValue", "value": "test_data/sample/pics", "name": "Directory"}, {"@type": "PropertyValue", "value": "Off", "name": "DistortionCorrection"}, {"@type": "PropertyValue", "value": "Off", "name": "DistortionCorrection2"}, {"@type": "PropertyValue", "value": "Single Shot", "name": "DriveMode"}, {"@type": "PropertyValue", "value": "Baseline DCT, Huffman coding", "name": "EncodingProcess"}, {"@type": "PropertyValue", "value": "0100", "name": "EquipmentVersion"}, {"@type": "PropertyValue", "value": "Little-endian (Intel, II)", "name": "ExifByteOrder"}, {"@type": "PropertyValue", "value": 2688, "name": "ExifImageHeight"}, {"@type": "PropertyValue", "value": 4032, "name": "ExifImageWidth"}, {"@type": "PropertyValue", "value": 10.8, "name": "ExifToolVersion"}, {"@type": "PropertyValue", "value": "0221", "name": "ExifVersion"}, {"@type": "PropertyValue", "value": -0.3, "name": "ExposureCompensation"}, {"@type": "PropertyValue", "value": "Manual", "name": "ExposureMode"}, {"@type": "PropertyValue", "value": "Aperture-priority AE", "name": "ExposureProgram"}, {"@type": "PropertyValue", "value": 0, "name": "ExposureShift"}, {"@type": "PropertyValue", "value": "1/80", "name": "ExposureTime"}, {"@type": "PropertyValue", "value": "Off", "name": "ExtendedWBDetect"}, {"@type": "PropertyValue", "value": "None", "name": "Extender"}, {"@type": "PropertyValue", "value": 0, "name": "ExtenderFirmwareVersion"}, {"@type": "PropertyValue", "value": "", "name": "ExtenderModel"}, {"@type": "PropertyValue", "value": "", "name": "ExtenderSerialNumber"}, {"@type": "PropertyValue", "value": "Not attached", "name": "ExtenderStatus"}, {"@type": "PropertyValue", "value": "Off", "name": "ExternalFlash"}, {"@type": "PropertyValue", "value": "Bounce or Off", "name": "ExternalFlashBounce"}, {"@type": "PropertyValue", "value": 0, "name": "ExternalFlashZoom"}, {"@type": "PropertyValue", "value": 6.3, "name": "FNumber"}, {"@type": "PropertyValue", "value": "47.8 deg (1.22 m)", "name": "FOV"}, {"@type": "PropertyValue", "value": "(Binary data 383 bytes, use -b option to extract)", "name": "FaceDetectArea"}, {"@type": "PropertyValue", "value": "0 14 320 212 0 14 320 212 0 0 0 0", "name": "FaceDetectFrameCrop"}, {"@type": "PropertyValue", "value": "320 240 320 240 0 0", "name": "FaceDetectFrameSize"}, {"@type": "PropertyValue", "value": "0 0 0", "name": "FacesDetected"}, {"@type": "PropertyValue", "value": "2018:09:19 16:05:34+10:00", "name": "FileAccessDate"}, {"@type": "PropertyValue", "value": "2018:09:19 16:04:32+10:00", "name": "FileInodeChangeDate"}, {"@type": "PropertyValue", "value": "2017:06:18 13:10:09+10:00", "name": "FileModifyDate"}, {"@type": "PropertyValue", "value": "2017-06-11 12.56.14.jpg", "name": "FileName"}, {"@type": "PropertyValue", "value": "rw-r--r--", "name": "FilePermissions"}, {"@type": "PropertyValue", "value": "4.9 MB", "name": "FileSize"}, {"@type": "PropertyValue", "value": "Digital Camera", "name": "FileSource"}, {"@type": "PropertyValue", "value": "JPEG", "name": "FileType"}, {"@type": "PropertyValue", "value": "jpg", "name": "FileTypeExtension"}, {"@type": "PropertyValue", "value": "On, Did not fire", "name": "Flash"}, {"@type": "PropertyValue", "value": "Off; 0; 0", "name": "FlashControlMode"}, {"@type": "PropertyValue", "value": 0, "name": "FlashExposureComp"}, {"@type": "PropertyValue", "value": 0, "name": "FlashFirmwareVersion"}, {"@type": "PropertyValue", "value": "n/a", "name": "FlashIntensity"}, {"@type": "PropertyValue", "value": "Fill-in", "name": "FlashMode"}, {"@type": "PropertyValue", "value": "None", "name": "FlashModel"}, {"@type": "PropertyValue", "value": "Off", "name": "FlashRemoteControl"}, {"@type": "PropertyValue", "value": "", "name": "FlashSerialNumber"}, {"@type": "PropertyValue", "value": "None", "name": "FlashType"}, {"@type": "PropertyValue", "value": "0100", "name": "FlashpixVersion"}, {"@type": "PropertyValue", "value": "20.0 mm", "name": "FocalLength"}, {"@type": "PropertyValue", "value": "20.0 mm (35 mm equivalent: 40.1 mm)", "name": "FocalLength35efl"}, {"@type": "PropertyValue", "value": "21.6 mm", "name": "FocalPlaneDiagonal"}, {"@type": "PropertyValue", "value": "1.375 m", "name": "FocusDistance"}, {"@type": "PropertyValue", "value": "0100", "name": "FocusInfoVersion"}, {"@type": "PropertyValue", "value": "Single AF; S-AF, Imager AF", "name": "FocusMode"}, {"@type": "PropertyValue", "value": "AF Used; 64", "name": "FocusProcess"}, {"@type": "PropertyValue", "value": 3151, "name": "FocusStepCount"}, {"@type": "PropertyValue", "value": 256, "name": "GainBase"}, {"@type": "PropertyValue", "value": "High gain up", "name": "GainControl"}, {"@type": "PropertyValue", "value": "Normal; Auto-Override", "name": "Gradation"}, {"@type": "PropertyValue", "value": "4.23 m", "name": "HyperfocalDistance"}, {"@type": "PropertyValue", "value": 400, "name": "ISO"}, {"@type": "PropertyValue", "value": "OLYMPUS DIGITAL CAMERA         ", "name": "ImageDescription"}, {"@type": "PropertyValue", "value": 2688, "name": "ImageHeight"}, {"@type": "PropertyValue", "value": "0112", "name": "ImageProcessingVersion"}, {"@type": "PropertyValue", "value": "4032x2688", "name": "ImageSize"}, {"@type": "PropertyValue", "value": "On, Mode 1", "name": "ImageStabilization"}, {"@type": "PropertyValue", "value": 4032, "name": "ImageWidth"}, {"@type": "PropertyValue", "value": "Off", "name": "InternalFlash"}, {"@type": "PropertyValue", "value": "4102011002108002               ", "name": "InternalSerialNumber"}, {"@type": "PropertyValue", "value": "R98 - DCF basic file (sRGB)", "name": "InteropIndex"}, {"@type": "PropertyValue", "value": "0100", "name": "InteropVersion"}, {"@type": "PropertyValue", "value": 1.1, "name": "LensFirmwareVersion"}, {"@type": "PropertyValue", "value": "Lumix G 20mm F1.7 Asph.", "name": "LensID"}, {"@type": "PropertyValue", "value": "LUMIX G 20/F1.7", "name": "LensModel"}, {"@type": "PropertyValue", "value": "0x4110", "name": "LensProperties"}, {"@type": "PropertyValue", "value": "01FG3033651", "name": "LensSerialNumber"}, {"@type": "PropertyValue", "value": "Lumix G 20mm F1.7 Asph.", "name": "LensType"}, {"@type": "PropertyValue", "value": "Unknown", "name": "LightSource"}, {"@type": "PropertyValue", "value": 9.6, "name": "LightValue"}, {"@type": "PropertyValue", "value": "image/jpeg", "name": "MIMEType"}, {"@type": "PropertyValue", "value": "Off", "name": "MacroLED"}, {"@type": "PropertyValue", "value": "Off", "name": "MacroMode"}, {"@type": "PropertyValue", "value": "OLYMPUS IMAGING CORP.", "name": "Make"}, {"@type": "PropertyValue", "value": "0 kPa", "name": "ManometerPressure"}, {"@type": "PropertyValue", "value": "0 m, 0 ft", "name": "ManometerReading"}, {"@type": "PropertyValue", "value": "Off", "name": "ManualFlash"}, {"@type": "PropertyValue", "value": "n/a", "name": "ManualFlashStrength"}, {"@type": "PropertyValue", "value": 1.8, "name": "MaxAperture"}, {"@type": "PropertyValue", "value": 1.7, "name": "MaxApertureAtMaxFocal"}, {"@type": "PropertyValue", "value": 1.7, "name": "MaxApertureAtMinFocal"}, {"@type": "PropertyValue", "value": 1.7, "name": "MaxApertureValue"}, {"@type": "PropertyValue", "value": "8 8 0", "name": "MaxFaces"}, {"@type": "PropertyValue", "value": 20, "name": "MaxFocalLength"}, {"@type": "PropertyValue", "value": 10.8, "name": "Megapixels"}, {"@type": "PropertyValue", "value": "ESP", "name": "MeteringMode"}, {"@type": "PropertyValue", "value": 20, "name": "MinFocalLength"}, {"@type": "PropertyValue", "value": "E-PL1", "name": "Model"}, {"@type": "PropertyValue", "value": "Off", "name": "ModifiedSaturation"}, {"@type": "PropertyValue", "value": "2017:06:11 12:56:14", "name": "ModifyDate"}, {"@type": "PropertyValue", "value": "Off; 1", "name": "MultipleExposureMode"}, {"@type": "PropertyValue", "value": "Standard", "name": "NoiseFilter"}, {"@type": "PropertyValue", "value": "Auto", "name": "NoiseReduction"}, {"@type": "PropertyValue", "value": "(none)", "name": "NoiseReduction2"}, {"@type": "PropertyValue", "value": "Horizontal (normal)", "name": "Orientation"}, {"@type": "PropertyValue", "value": "Off", "name": "PanoramaMode"}, {"@type": "PropertyValue", "value": "Vivid; 2", "name": "PictureMode"}, {"@type": "PropertyValue", "value": "n/a", "name": "PictureModeBWFilter"}, {"@type": "PropertyValue", "value": "0 (min -2, max 2)", "name": "PictureModeContrast"}, {"@type": "PropertyValue", "value": "Standard", "name": "PictureModeEffect"}, {"@type": "PropertyValue", "value": "0 (min -2, max 2)", "name": "PictureModeSaturation"}, {"@type": "PropertyValue", "value": "0 (min -2, max 2)", "name": "PictureModeSharpness"}, {"@type": "PropertyValue", "value": "n/a", "name": "PictureModeTone"}, {"@type": "PropertyValue", "value": "(Binary data 514170 bytes, use -b option to extract)", "name": "PreviewImage"}, {"@type": "PropertyValue", "value": 514170, "name": "PreviewImageLength"}, {"@type": "PropertyValue", "value": 4600608, "name": "PreviewImageStart"}, {"@type": "PropertyValue", "value": "Yes", "name": "PreviewImageValid"}, {"@type": "PropertyValue", "value": "0300", "name": "PrintIMVersion"}, {"@type": "PropertyValue", "value": "sRGB", "name": "RawDevColorSpace"}, {"@type": "PropertyValue", "value": "0 0 0", "name": "RawDevContrastValue"}, {"@type": "PropertyValue", "value": "Original", "name": "RawDevEditStatus"}, {"@type": "PropertyValue", "value": "High Speed", "name": "RawDevEngine"}, {"@type": "PropertyValue", "value": 0, "name": "RawDevExposureBiasValue"}, {"@type": "PropertyValue", "value": "0 0 0", "name": "RawDevGrayPoint"}, {"@type": "PropertyValue", "value": 0, "name": "RawDevMemoryColorEmphasis"}, {"@type": "PropertyValue", "value": "(none)", "name": "RawDevNoiseReduction"}, {"@type": "PropertyValue", "value": "0 0 0", "name": "RawDevSaturationEmphasis"}, {"@type": "PropertyValue", "value": "(none)", "name": "RawDevSettings"}, {"@type": "PropertyValue", "value": "0 0 0", "name": "RawDevSharpnessValue"}, {"@type": "PropertyValue", "value": "0100", "name": "RawDevVersion"}, {"@type": "PropertyValue", "value": 0, "name": "RawDevWBFineAdjustment"}, {"@type": "PropertyValue", "value": 0, "name": "RawDevWhiteBalanceValue"}, {"@type": "PropertyValue", "value": 1.78125, "name": "RedBalance"}, {"@type": "PropertyValue", "value": "inches", "name": "ResolutionUnit"}, {"@type": "PropertyValue", "value": "High", "name": "Saturation"}, {"@type": "PropertyValue", "value": 2.0, "name": "ScaleFactor35efl"}, {"@type": "PropertyValue", "value": "Standard", "name": "SceneCaptureType"}, {"@type": "PropertyValue", "value": 0, "name": "SceneDetect"}, {"@type": "PropertyValue", "value": "Standard", "name": "SceneMode"}, {"@type": "PropertyValue", "value": "0 0", "name": "SensorCalibration"}, {"@type": "PropertyValue", "value": "48.0 C", "name": "SensorTemperature"}, {"@type": "PropertyValue", "value": "B3B505338", "name": "SerialNumber"}, {"@type": "PropertyValue", "value": "Off", "name": "ShadingCompensation"}, {"@type": "PropertyValue", "value": "Off", "name": "ShadingCompensation2"}, {"@type": "PropertyValue", "value": "Hard", "name": "Sharpness"}, {"@type": "PropertyValue", "value": "1 (min -5, max 5)", "name": "SharpnessSetting"}, {"@type": "PropertyValue", "value": "1/80", "name": "ShutterSpeed"}, {"@type": "PropertyValue", "value": "Version 1.1", "name": "Software"}, {"@type": "PropertyValue", "value": "test_data/sample/pics/2017-06-11 12.56.14.jpg", "name": "SourceFile"}, {"@type": "PropertyValue", "value": "Normal, Sequence: 0, Panorama: (none)", "name": "SpecialMode"}, {"@type": "PropertyValue", "value": "(Binary data 7738 bytes, use -b option to extract)", "name": "ThumbnailImage"}, {"@type": "PropertyValue", "value": 7738, "name": "ThumbnailLength"}, {"@type": "PropertyValue", "value": 14464, "name": "ThumbnailOffset"}, {"@type": "PropertyValue", "value": "0; 0; 0; 0; 0; 0; 0; 0; 0; 0; 0; 0", "name": "ToneLevel"}, {"@type": "PropertyValue", "value": "", "name": "UserComment"}, {"@type": "PropertyValue", "value": "456 350 256 260", "name": "WB_RBLevels"}, {"@type": "PropertyValue", "value": "Auto", "name": "WhiteBalance"}, {"@type": "PropertyValue", "value": "Auto", "name": "WhiteBalance2"}, {"@type": "PropertyValue", "value": "0 0", "name": "WhiteBalanceBracket"}, {"@type": "PropertyValue", "value": "Auto", "name": "WhiteBalanceTemperature"}, {"@type": "PropertyValue", "value": 314, "name": "XResolution"}, {"@type": "PropertyValue", "value": "Co-sited", "name": "YCbCrPositioning"}, {"@type": "PropertyValue", "value": "YCbCr4:2:2 (2 1)", "name": "YCbCrSubSampling"}, {"@type": "PropertyValue", "value": 314, "name": "YResolution"}, {"@type": "PropertyValue", "value": 0, "name": "ZoomStepCount"}]'

# --- Row 3: brand new record for the sepia version of the image ---
$ws.Range("A3").Value = "sepia_fence.jpg"
$ws.Range("B3").Value = "pics/thumbs/sepia_fence.png"
$ws.Range("C3").Value = 'Sepia tone version of my fence/dog pic'
$ws.Range("G3").Value = "ImageObject"
$ws.Range("H3").Value = '[{"@type": "PropertyValue", "value": "Off", "name": "AELock"}, {"@type": "PropertyValue", "value": "(111,106)-(144,149)", "name": "AFAreas"}, {"@type": "PropertyValue", "value": "Off", "name": "AFFineTune"}, {"@type": "PropertyValue", "value": "0 0 0", "name": "AFFineTuneAdj"}, {"@type": "PropertyValue", "value": "Left (or n/a)", "name": "AFPoint"}, {"@type": "PropertyValue", "value": "(49%,69%) (49%,69%)", "name": "AFPointSelected"}, {"@type": "PropertyValue", "value": "Ready", "name": "AFSearch"}, {"@type": "PropertyValue", "value": 6.3, "name": "Aperture"}, {"@type": "PropertyValue", "value": "Off; 0; 0; 0", "name": "ArtFilter"}, {"@type": "PropertyValue", "value": "0 0 4031 2687", "name": "AspectFrame"}, {"@type": "PropertyValue", "value": "3:2", "name": "AspectRatio"}, {"@type": "PropertyValue", "value": 8, "name": "BitsPerSample"}, {"@type": "PropertyValue", "value": "64 64 64 64", "name": "BlackLevel2"}, {"@type": "PropertyValue", "value": 1.367188, "name": "BlueBalance"}, {"@type": "PropertyValue", "value": 1.103, "name": "BodyFirmwareVersion"}, {"@type": "PropertyValue", "value": "OLYMPUS DIGITAL CAMERA         ", "name": "CameraID"}, {"@type": "PropertyValue", "value": "0100", "name": "CameraSettingsVersion"}, {"@type": "PropertyValue", "value": "E-PL1", "name": "CameraType2"}, {"@type": "PropertyValue", "value": "0.015 mm", "name": "CircleOfConfusion"}, {"@type": "PropertyValue", "value": 3, "name": "ColorComponents"}, {"@type": "PropertyValue", "value": "322 -40 -26 -32 332 -44 -2 -72 330", "name": "ColorMatrix"}, {"@type": "PropertyValue", "value": "sRGB", "name": "ColorSpace"}, {"@type": "PropertyValue", "value": "Y, Cb, Cr, -", "name": "ComponentsConfiguration"}, {"@type": "PropertyValue", "value": "JPEG (old-style)", "name": "Compression"}, {"@type": "PropertyValue", "value": 4, "name": "CompressionFactor"}, {"@type": "PropertyValue", "value": "Normal", "name": "Contrast"}, {"@type": "PropertyValue", "value": "0 (min -5, max 5)", "name": "ContrastSetting"}, {"@type": "PropertyValue", "value": "2017:06:11 12:56:14", "name": "CreateDate"}, {"@type": "PropertyValue", "value": 3024, "name": "CropHeight"}, {"@type": "PropertyValue", "value": "28 0", "name": "CropLeft"}, {"@type": "PropertyValue", "value": "30 0", "name": "CropTop"}, {"@type": "PropertyValue", "value": 4032, "name": "CropWidth"}, {"@type": "PropertyValue", "value": "Normal", "name": "CustomRendered"}, {"@type": "PropertyValue", "value": "1 (min -5, max 5)", "name": "CustomSaturation"}, {"@type": "PropertyValue", "value": "0.98 m (1.04 - 2.02 m)", "name": "DOF"}, {"@type": "PropertyValue", "value": "2017:06:11 12:56:14", "name": "DateTimeOriginal"}, {"@type": "PropertyValue", "value": 1, "name": "DigitalZoomRatio"}, {"@type": "PropertyValue", "value": "test_data/sample/pics", "name": "Directory"}, {"@type": "PropertyValue", "value": "Off", "name": "DistortionCorrection"}, {"@type": "PropertyValue", "value": "Off", "name": "DistortionCorrection2"}, {"@type": "PropertyValue", "value": "Single Shot", "name": "DriveMode"}, {"@type": "PropertyValue", "value": "Baseline DCT, Huffman coding", "name": "EncodingProcess"}, {"@type": "PropertyValue", "value": "0100", "name": "EquipmentVersion"}, {"@type": "PropertyValue", "value": "Little-endian (Intel, II)", "name": "ExifByteOrder"}, {"@type": "PropertyValue", "value": 2688, "name": "ExifImageHeight"}, {"@type": "PropertyValue", "value": 4032, "name": "ExifImageWidth"}, {"@type": "PropertyValue", "value": 10.8, "name": "ExifToolVersion"}, {"@type": "PropertyValue", "value": "0221", "name": "ExifVersion"}, {"@type": "PropertyValue", "value": -0.3, "name": "ExposureCompensation"}, {"@type": "PropertyValue", "value": "Manual", "name": "ExposureMode"}, {"@type": "PropertyValue", "value": "Aperture-priority AE", "name": "ExposureProgram"}, {"@type": "PropertyValue", "value": 0, "name": "ExposureShift"}, {"@type": "PropertyValue", "value": "1/80", "name": "ExposureTime"}, {"@type": "PropertyValue", "value": "Off", "name": "ExtendedWBDetect"}, {"@type": "PropertyValue", "value": "None", "name": "Extender"}, {"@type": "PropertyValue", "value": 0, "name": "ExtenderFirmwareVersion"}, {"@type": "PropertyValue", "value": "", "name": "ExtenderModel"}, {"@type": "PropertyValue", "value": "", "name": "ExtenderSerialNumber"}, {"@type": "PropertyValue", "value": "Not attached", "name": "ExtenderStatus"}, {"@type": "PropertyValue", "value": "Off", "name": "ExternalFlash"}, {"@type": "PropertyValue", "value": "Bounce or Off", "name": "ExternalFlashBounce"}, {"@type": "PropertyValue", "value": 0, "name": "ExternalFlashZoom"}, {"@type": "PropertyValue", "value": 6.3, "name": "FNumber"}, {"@type": "PropertyValue", "value": "47.8 deg (1.22 m)", "name": "FOV"}, {"@type": "PropertyValue", "value": "(Binary data 383 bytes, use -b option to extract)", "name": "FaceDetectArea"}, {"@type": "PropertyValue", "value": "0 14 320 212 0 14 320 212 0 0 0 0", "name": "FaceDetectFrameCrop"}, {"@type": "PropertyValue", "value": "320 240 320 240 0 0", "name": "FaceDetectFrameSize"}, {"@type": "PropertyValue", "value": "0 0 0", "name": "FacesDetected"}, {"@type": "PropertyValue", "value": "2018:09:19 17:01:08+10:00", "name": "FileAccessDate"}, {"@type": "PropertyValue", "value": "2018:09:19 17:01:07+10:00", "name": "FileInodeChangeDate"}, {"@type": "PropertyValue", "value": "2018:09:19 17:01:07+10:00", "name": "FileModifyDate"}, {"@type": "PropertyValue", "value": "sepia_fence.jpg", "name": "FileName"}, {"@type": "PropertyValue", "value": "rw-r--r--", "name": "FilePermissions"}, {"@type": "PropertyValue", "value": "4.6 MB", "name": "FileSize"}, {"@type": "PropertyValue", "value": "Digital Camera", "name": "FileSource"}, {"@type": "PropertyValue", "value": "JPEG", "name": "FileType"}, {"@type": "PropertyValue", "value": "jpg", "name": "FileTypeExtension"}, {"@type": "PropertyValue", "value": "On, Did not fire", "name": "Flash"}, {"@type": "PropertyValue", "value": "Off; 0; 0", "name": "FlashControlMode"}, {"@type": "PropertyValue", "value": 0, "name": "FlashExposureComp"}, {"@type": "PropertyValue", "value": 0, "name": "FlashFirmwareVersion"}, {"@type": "PropertyValue", "value": "n/a", "name": "FlashIntensity"}, {"@type": "PropertyValue", "value": "Fill-in", "name": "FlashMode"}, {"@type": "PropertyValue", "value": "None", "name": "FlashModel"}, {"@type": "PropertyValue", "value": "Off", "name": "FlashRemoteControl"}, {"@type": "PropertyValue", "value": "", "name": "FlashSerialNumber"}, {"@type": "PropertyValue", "value": "None", "name": "FlashType"}, {"@type": "PropertyValue", "value": "0100", "name": "FlashpixVersion"}, {"@type": "PropertyValue", "value": "20.0 mm", "name": "FocalLength"}, {"@type": "PropertyValue", "value": "20.0 mm (35 mm equivalent: 40.1 mm)", "name": "FocalLength35efl"}, {"@type": "PropertyValue", "value": "21.6 mm", "name": "FocalPlaneDiagonal"}, {"@type": "PropertyValue", "value": "1.375 m", "name": "FocusDistance"}, {"@type": "PropertyValue", "value": "0100", "name": "FocusInfoVersion"}, {"@type": "PropertyValue", "value": "Single AF; S-AF, Imager AF", "name": "FocusMode"}, {"@type": "PropertyValue", "value": "AF Used; 64", "name": "FocusProcess"}, {"@type": "PropertyValue", "value": 3151, "name": "FocusStepCount"}, {"@type": "PropertyValue", "value": 256, "name": "GainBase"}, {"@type": "PropertyValue", "value": "High gain up", "name": "GainControl"}, {"@type": "PropertyValue", "value": "Normal; Auto-Override", "name": "Gradation"}, {"@type": "PropertyValue", "value": "4.23 m", "name": "HyperfocalDistance"}, {"@type": "PropertyValue", "value": 400, "name": "ISO"}, {"@type": "PropertyValue", "value": "OLYMPUS DIGITAL CAMERA         ", "name": "ImageDescription"}, {"@type": "PropertyValue", "value": 2688, "name": "ImageHeight"}, {"@type": "PropertyValue", "value": "0112", "name": "ImageProcessingVersion"}, {"@type": "PropertyValue", "value": "4032x2688", "name": "ImageSize"}, {"@type": "PropertyValue", "value": "On, Mode 1", "name": "ImageStabilization"}, {"@type": "PropertyValue", "value": 4032, "name": "ImageWidth"}, {"@type": "PropertyValue", "value": "Off", "name": "InternalFlash"}, {"@type": "PropertyValue", "value": "4102011002108002               ", "name": "InternalSerialNumber"}, {"@type": "PropertyValue", "value": "R98 - DCF basic file (sRGB)", "name": "InteropIndex"}, {"@type": "PropertyValue", "value": "0100", "name": "InteropVersion"}, {"@type": "PropertyValue", "value": 1.01, "name": "JFIFVersion"}, {"@type": "PropertyValue", "value": 1.1, "name": "LensFirmwareVersion"}, {"@type": "PropertyValue", "value": "Lumix G 20mm F1.7 Asph.", "name": "LensID"}, {"@type": "PropertyValue", "value": "LUMIX G 20/F1.7", "name": "LensModel"}, {"@type": "PropertyValue", "value": "0x4110", "name": "LensProperties"}, {"@type": "PropertyValue", "value": "01FG3033651", "name": "LensSerialNumber"}, {"@type": "PropertyValue", "value": "Lumix G 20mm F1.7 Asph.", "name": "LensType"}, {"@type": "PropertyValue", "value": "Unknown", "name": "LightSource"}, {"@type": "PropertyValue", "value": 9.6, "name": "LightValue"}, {"@type": "PropertyValue", "value": "image/jpeg", "name": "MIMEType"}, {"@type": "PropertyValue", "value": "Off", "name": "MacroLED"}, {"@type": "PropertyValue", "value": "Off", "name": "MacroMode"}, {"@type": "PropertyValue", "value": "OLYMPUS IMAGING CORP.", "name": "Make"}, {"@type": "PropertyValue", "value": "0 kPa", "name": "ManometerPressure"}, {"@type": "PropertyValue", "value": "0 m, 0 ft", "name": "ManometerReading"}, {"@type": "PropertyValue", "value": "Off", "name": "ManualFlash"}, {"@type": "PropertyValue", "value": "n/a", "name": "ManualFlashStrength"}, {"@type": "PropertyValue", "value": 1.8, "name": "MaxAperture"}, {"@type": "PropertyValue", "value": 1.7, "name": "MaxApertureAtMaxFocal"}, {"@type": "PropertyValue", "value": 1.7, "name": "MaxApertureAtMinFocal"}, {"@type": "PropertyValue", "value": 1.7, "name": "MaxApertureValue"}, {"@type": "PropertyValue", "value": "8 8 0", "name": "MaxFaces"}, {"@type": "PropertyValue", "value": 20, "name": "MaxFocalLength"}, {"@type": "PropertyValue", "value": 10.8, "name": "Megapixels"}, {"@type": "PropertyValue", "value": "ESP", "name": "MeteringMode"}, {"@type": "PropertyValue", "value": 20, "name": "MinFocalLength"}, {"@type": "PropertyValue", "value": "E-PL1", "name": "Model"}, {"@type": "PropertyValue", "value": "Off", "name": "ModifiedSaturation"}, {"@type": "PropertyValue", "value": "2017:06:11 12:56:14", "name": "ModifyDate"}, {"@type": "PropertyValue", "value": "Off; 1", "name": "MultipleExposureMode"}, {"@type": "PropertyValue", "value": "Standard", "name": "NoiseFilter"}, {"@type": "PropertyValue", "value": "Auto", "name": "NoiseReduction"}, {"@type": "PropertyValue", "value": "(none)", "name": "NoiseReduction2"}, {"@type": "PropertyValue", "value": "Horizontal (normal)", "name": "Orientation"}, {"@type": "PropertyValue", "value": "Off", "name": "PanoramaMode"}, {"@type": "PropertyValue", "value": "Vivid; 2", "name": "PictureMode"}, {"@type": "PropertyValue", "value": "n/a", "name": "PictureModeBWFilter"}, {"@type": "PropertyValue", "value": "0 (min -2, max 2)", "name": "PictureModeContrast"}, {"@type": "PropertyValue", "value": "Standard", "name": "PictureModeEffect"}, {"@type": "PropertyValue", "value": "0 (min -2, max 2)", "name": "PictureModeSaturation"}, {"@type": "PropertyValue", "value": "0 (min -2, max 2)", "name": "PictureModeSharpness"}, {"@type": "PropertyValue", "value": "n/a", "name": "PictureModeTone"}, {"@type": "PropertyValue", "value": "(Binary data 514170 bytes, use -b option to extract)", "name": "PreviewImage"}, {"@type": "PropertyValue", "value": 514170, "name": "PreviewImageLength"}, {"@type": "PropertyValue", "value": 4600626, "name": "PreviewImageStart"}, {"@type": "PropertyValue", "value": "Yes", "name": "PreviewImageValid"}, {"@type": "PropertyValue", "value": "0300", "name": "PrintIMVersion"}, {"@type": "PropertyValue", "value": "sRGB", "name": "RawDevColorSpace"}, {"@type": "PropertyValue", "value": "0 0 0", "name": "RawDevContrastValue"}, {"@type": "PropertyValue", "value": "Original", "name": "RawDevEditStatus"}, {"@type": "PropertyValue", "value": "High Speed", "name": "RawDevEngine"}, {"@type": "PropertyValue", "value": 0, "name": "RawDevExposureBiasValue"}, {"@type": "PropertyValue", "value": "0 0 0", "name": "RawDevGrayPoint"}, {"@type": "PropertyValue", "value": 0, "name": "RawDevMemoryColorEmphasis"}, {"@type": "PropertyValue", "value": "(none)", "name": "RawDevNoiseReduction"}, {"@type": "PropertyValue", "value": "0 0 0", "name": "RawDevSaturationEmphasis"}, {"@type": "PropertyValue", "value": "(none)", "name": "RawDevSettings"}, {"@type": "PropertyValue", "value": "0 0 0", "name": "RawDevSharpnessValue"}, {"@type": "PropertyValue", "value": "0100", "name": "RawDevVersion"}, {"@type": "PropertyValue", "value": 0, "name": "RawDevWBFineAdjustment"}, {"@type": "PropertyValue", "value": 0, "name": "RawDevWhiteBalanceValue"}, {"@type": "PropertyValue", "value": 1.78125, "name": "RedBalance"}, {"@type": "PropertyValue", "value": "inches", "name": "ResolutionUnit"}, {"@type": "PropertyValue", "value": "High", "name": "Saturation"}, {"@type": "PropertyValue", "value": 2.0, "name": "ScaleFactor35efl"}, {"@type": "PropertyValue", "value": "Standard", "name": "SceneCaptureType"}, {"@type": "PropertyValue", "value": 0, "name": "SceneDetect"}, {"@type": "PropertyValue", "value": "Standard", "name": "SceneMode"}, {"@type": "PropertyValue", "value": "0 0", "name": "SensorCalibration"}, {"@type": "PropertyValue", "value": "48.0 C", "name": "SensorTemperature"}, {"@type": "PropertyValue", "value": "B3B505338", "name": "SerialNumber"}, {"@type": "PropertyValue", "value": "Off", "name": "ShadingCompensation"}, {"@type": "PropertyValue", "value": "Off", "name": "ShadingCompensation2"}, {"@type": "PropertyValue", "value": "Hard", "name": "Sharpness"}, {"@type": "PropertyValue", "value": "1 (min -5, max 5)", "name": "SharpnessSetting"}, {"@type": "PropertyValue", "value": "1/80", "name": "ShutterSpeed"}, {"@type": "PropertyValue", "value": "Version 1.1", "name": "Software"}, {"@type": "PropertyValue", "value": "test_data/sample/pics/sepia_fence.jpg", "name": "SourceFile"}, {"@type": "PropertyValue", "value": "Normal, Sequence: 0, Panorama: (none)", "name": "SpecialMode"}, {"@type": "PropertyValue", "value": "(Binary data 7738 bytes, use -b option to extract)", "name": "ThumbnailImage"}, {"@type": "PropertyValue", "value": 7738, "name": "ThumbnailLength"}, {"@type": "PropertyValue", "value": 14482, "name": "ThumbnailOffset"}, {"@type": "PropertyValue", "value": "0; 0; 0; 0; 0; 0; 0; 0; 0; 0; 0; 0", "name": "ToneLevel"}, {"@type": "PropertyValue", "value": "", "name": "UserComment"}, {"@type": "PropertyValue", "value": "456 350 256 260", "name": "WB_RBLevels"}, {"@type": "PropertyValue", "value": "Auto", "name": "WhiteBalance"}, {"@type": "PropertyValue", "value": "Auto", "name": "WhiteBalance2"}, {"@type": "PropertyValue", "value": "0 0", "name": "WhiteBalanceBracket"}, {"@type": "PropertyValue", "value": "Auto", "name": "WhiteBalanceTemperature"}, {"@type": "PropertyValue", "value": 314, "name": "XResolution"}, {"@type": "PropertyValue", "value": "Co-sited", "name": "YCbCrPositioning"}, {"@type": "PropertyValue", "value": "YCbCr4:2:2 (2 1)", "name": "YCbCrSubSampling"}, {"@type": "PropertyValue", "value": 314, "name": "YResolution"}, {"@type": "PropertyValue", "value": 0, "name": "ZoomStepCount"}]'

